# Refresh cryptos list: updated Price (D) and Volume(1h) (E) columns
# with the latest scrape values (GitHub Actions crypto-list updater).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column prices are plain text in this sheet (e.g. "354.00", "67.097.39"),
# so a leading apostrophe is used to force text storage and avoid Excel
# auto-converting them to numbers (which would drop formatting like trailing
# zeros or multi-dot thousand separators).
$ws.Range("D2").Value = "'67.097.39"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "'2.544.74"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("D5").Value = "'590.16"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'173.33"
$ws.Range("E6").Value = "  +4.92%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D9").Value = "'2.544.46"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "'0.346"
$ws.Range("E13").Value = "  -5.15%  "
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "'3.011.57"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "'66.988.16"
$ws.Range("D18").Value = "'2.533.38"
$ws.Range("E18").Value = "  -2.13%  "
$ws.Range("D19").Value = "'8.03"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").Value = "'11.29"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Value = "'354.00"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  +6.09%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'69.76"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "'10.09"
$ws.Range("E27").Value = "  -3.78%  "
$ws.Range("D28").Value = "'2.678.58"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'533.43"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'157.18"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").Value = "'18.66"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").Value = "'18.45"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("D46").Value = "'39.72"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "'149.59"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").Value = "'0.0₆0277"
$ws.Range("E49").Value = "  -4.91%  "
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("E51").Value = "  +0.30%  "
